$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append per-play yardage logs (Week 15 actual + Week 16 sim)
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 3 4 0 0 4 10 1 1 12 4 5 0 8 2 5 4 1 6 6 9 3 7 2 0 11 3 2 0 3 24 8 30 4 12 1 4"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 5 1 9 37 15 12 4 27 20 20 7 1 17 9 11 1"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 6 8 11 1 2 2 3 1 5 3 5 7 4 1 -5 11 0 4 6 6 20 -2 2 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 16 2 14 19 4 11 6 23 1 9 24 10"

# ---------------------------------------------------------------------------
# OFF sheet: updated situational / play totals
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 338
$offWs.Range("D2").Value = 24
$offWs.Range("F2").Value = 117
$offWs.Range("G2").Value = 106
$offWs.Range("J2").Value = 62
$offWs.Range("L2").Value = 575
$offWs.Range("M2").Value = 378
$offWs.Range("O2").Value = 28
$offWs.Range("P2").Value = 19
$offWs.Range("Q2").Value = 1042

$offWs.Range("C3").Value = 377
$offWs.Range("D3").Value = 13
$offWs.Range("E3").Value = 47
$offWs.Range("F3").Value = 201
$offWs.Range("I3").Value = 98
$offWs.Range("J3").Value = 123
$offWs.Range("N3").Value = 36

# ---------------------------------------------------------------------------
# DEF sheet: updated situational / play totals
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 360
$defWs.Range("D2").Value = 19
$defWs.Range("F2").Value = 101
$defWs.Range("G2").Value = 87
$defWs.Range("H2").Value = 8
$defWs.Range("J2").Value = 40
$defWs.Range("L2").Value = 577
$defWs.Range("M2").Value = 343
$defWs.Range("O2").Value = 42
$defWs.Range("Q2").Value = 991

$defWs.Range("B3").Value = 15
$defWs.Range("C3").Value = 364
$defWs.Range("E3").Value = 68
$defWs.Range("F3").Value = 218
$defWs.Range("G3").Value = 64
$defWs.Range("H3").Value = 62
$defWs.Range("I3").Value = 113
$defWs.Range("N3").Value = 50

# ---------------------------------------------------------------------------
# ST sheet: kickoff / punt counters + logs
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 150
$stWs.Range("D2").Value = 130
$stWs.Range("F2").Value = 132
$stWs.Range("G2").Value = 129
$stWs.Range("J2").Value = 54
$stWs.Range("K2").Value = 51
$stWs.Range("B3").Value = 110

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 51 52 56"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 11 12"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 0 0 0 5"

# ---------------------------------------------------------------------------
# TURNS sheet: interceptions / fumbles
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 14
$turnsWs.Range("D2").Value = 16
$turnsWs.Range("E2").Value = 16

# ---------------------------------------------------------------------------
# PEN sheet: penalties
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 32
$penWs.Range("B3").Value = 28
